$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BF column (58) holds a "Date" column that was stored as the text
# "5-7-2007-08" (day-month swapped / wrong year range) for every data row.
# Correct it to the real ISO date string "2008-05-07" for rows 2-31,
# keeping the value as plain text (not an auto-converted date serial).
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)
    $cell.NumberFormat = "@"
    $cell.Value = "2008-05-07"
    $cell.Style = "Normal"
}
